$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - update the ApplicationUrl value
$ws.Range("B2").Value = "https://demo.guru99.com/test/newtours/"

# Row 3 - becomes "UserName" / "DemoUser1"
$ws.Range("A3").Value = "UserName"
$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("A3").VerticalAlignment = -4108

$ws.Range("B3").ClearFormats()
$ws.Range("B3").Value = "DemoUser1"

# Row 4 - becomes "Password" / base64 value (copy A3's new alignment/font
# style instead of re-deriving it, so the style table doesn't pick up a
# stray intermediate alignment-only xf)
$ws.Range("A4").Value = "Password"
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B4").ClearFormats()
$ws.Range("B4").Value = "VGVzdFVzZXJAMTIzNDU="

# Rows 5-7 - the old wait.* labels and values are removed, only the B-column
# formatting remains (blank placeholder cells keep style 1)
$ws.Range("A5").ClearContents()
$ws.Range("B5").ClearContents()

$ws.Range("A6").ClearContents()
$ws.Range("B6").ClearContents()

$ws.Range("A7").ClearContents()
$ws.Range("B7").ClearContents()

# Update the active selection to A4, matching the saved workbook view
$ws.Range("A4").Select() | Out-Null
